$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 254.5
$ws.Range("I8").Value = 254.5
$ws.Range("K8").Value = 763.5
$ws.Range("M8").Value = -624.5

$ws.Range("H39").Value = 4295.222
$ws.Range("I39").Value = 1635.8
$ws.Range("J39").Value = 7619.5
$ws.Range("K39").Value = 4907.4
$ws.Range("L39").Value = 22858.5
$ws.Range("M39").Value = -4611.4
$ws.Range("N39").Value = -23450.5

$ws.Range("H105").Value = 19000
$ws.Range("J105").Value = 19000
$ws.Range("L105").Value = 19000
$ws.Range("N105").Value = -25988

$ws.Range("H135").Value = 2509.818
$ws.Range("I135").Value = 2509.818
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 22588.362
$ws.Range("L135").Value = 0
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -20053.362

$ws.Range("H138").Value = 2856.1355
$ws.Range("J138").Value = 3419.6047
$ws.Range("L138").Value = 10258.8141
$ws.Range("N138").Value = -20538.8141

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 340802.66
$ws.Range("I8").Value = 670336.7
$ws.Range("K8").Value = 670336.7
$ws.Range("M8").Value = -670192.7

$ws.Range("H61").Value = 6063
$ws.Range("I61").Value = 4704.684
$ws.Range("J61").Value = 14665.667
$ws.Range("K61").Value = 4704.684
$ws.Range("L61").Value = 14665.667
$ws.Range("M61").Value = -4492.684
$ws.Range("N61").Value = -15089.667

$ws.Range("H102").Value = 39079.555
$ws.Range("J102").Value = 100000
$ws.Range("L102").Value = 100000
$ws.Range("N102").Value = -103244

$ws.Range("H122").Value = 4755.085
$ws.Range("I122").Value = 4277.425
$ws.Range("K122").Value = 12832.275
$ws.Range("M122").Value = -10382.275

$ws.Range("H136").Value = 6063
$ws.Range("I136").Value = 4704.684
$ws.Range("J136").Value = 14665.667
$ws.Range("K136").Value = 14114.052
$ws.Range("L136").Value = 43997.001
$ws.Range("M136").Value = -11564.052
$ws.Range("N136").Value = -49097.001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1252.8422
$ws.Range("I86").Value = 844.0625
$ws.Range("J86").Value = 3433
$ws.Range("K86").Value = 844.0625
$ws.Range("L86").Value = 3433
$ws.Range("M86").Value = 278.9375
$ws.Range("N86").Value = -5679

$ws.Range("H89").Value = 1252.8422
$ws.Range("I89").Value = 844.0625
$ws.Range("J89").Value = 3433
$ws.Range("K89").Value = 4220.3125
$ws.Range("L89").Value = 17165
$ws.Range("M89").Value = 1395.6875
$ws.Range("N89").Value = -28397

$ws.Range("H107").Value = 2928.3572
$ws.Range("I107").Value = 3110.6667
$ws.Range("K107").Value = 3110.6667
$ws.Range("M107").Value = -1190.6667

$ws.Range("H134").Value = 3159.8696
$ws.Range("I134").Value = 3159.8696
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 9479.6088
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -6944.6088

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5584.8184
$ws.Range("I58").Value = 4929.125
$ws.Range("J58").Value = 7333.3335
$ws.Range("K58").Value = 4929.125
$ws.Range("L58").Value = 7333.3335
$ws.Range("M58").Value = -4726.125
$ws.Range("N58").Value = -7739.3335

$ws.Range("H110").Value = 49000.5
$ws.Range("J110").Value = 49000.5
$ws.Range("L110").Value = 49000.5
$ws.Range("N110").Value = -57180.5

$ws.Range("H122").Value = 103621.63
$ws.Range("I122").Value = 140111.31
$ws.Range("K122").Value = 420333.93
$ws.Range("M122").Value = -417883.93

$ws.Range("H136").Value = 5584.8184
$ws.Range("I136").Value = 4929.125
$ws.Range("J136").Value = 7333.3335
$ws.Range("K136").Value = 14787.375
$ws.Range("L136").Value = 22000.0005
$ws.Range("M136").Value = -12237.375
$ws.Range("N136").Value = -27100.0005

$ws.Range("H141").Value = 193141.72
$ws.Range("J141").Value = 193141.72
$ws.Range("L141").Value = 193141.72
$ws.Range("N141").Value = -203501.72

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 106900040
$ws.Range("I4").Value = 167983490
$ws.Range("J4").Value = 15274874
$ws.Range("K4").Value = 503950470
$ws.Range("L4").Value = 45824622
$ws.Range("M4").Value = -503950358
$ws.Range("N4").Value = -45824846

$ws.Range("H34").Value = 5550.5
$ws.Range("J34").Value = 8672.4
$ws.Range("L34").Value = 26017.2
$ws.Range("N34").Value = -26185.2

$ws.Range("H36").Value = 1293
$ws.Range("I36").Value = 1293
$ws.Range("K36").Value = 3879
$ws.Range("M36").Value = -3710

$ws.Range("H107").Value = 527.7778
$ws.Range("J107").Value = 533
$ws.Range("L107").Value = 1599
$ws.Range("N107").Value = -5439

$ws.Range("H122").Value = 2644.5264
$ws.Range("J122").Value = 2752.5557
$ws.Range("L122").Value = 24773.0013
$ws.Range("N122").Value = -29673.0013

$ws.Range("H139").Value = 3209.6191
$ws.Range("I139").Value = 3100.1052
$ws.Range("K139").Value = 9300.3156
$ws.Range("M139").Value = -4160.3156

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 141.14285
$ws.Range("I2").Value = 48
$ws.Range("J2").Value = 211
$ws.Range("K2").Value = 48
$ws.Range("L2").Value = 211
$ws.Range("M2").Value = 65
$ws.Range("N2").Value = -437

$ws.Range("H43").Value = 1499.5
$ws.Range("I43").Value = 1499.5
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 1499.5
$ws.Range("L43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -1348.5

$ws.Range("H102").Value = 2111.0386
$ws.Range("I102").Value = 1613.7142
$ws.Range("K102").Value = 1613.7142
$ws.Range("M102").Value = 8.285800000000108

$ws.Range("H122").Value = 7072.7856
$ws.Range("I122").Value = 6789.222
$ws.Range("J122").Value = 7583.2
$ws.Range("K122").Value = 20367.666
$ws.Range("L122").Value = 22749.6
$ws.Range("M122").Value = -17917.666
$ws.Range("N122").Value = -27649.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3666.3333
$ws.Range("I7").Value = 3249.5
$ws.Range("K7").Value = 3249.5
$ws.Range("M7").Value = -3137.5

$ws.Range("H40").Value = 7174.96
$ws.Range("I40").Value = 7371.5713
$ws.Range("J40").Value = 6142.75
$ws.Range("K40").Value = 7371.5713
$ws.Range("L40").Value = 6142.75
$ws.Range("M40").Value = -7235.5713
$ws.Range("N40").Value = -6414.75

$ws.Range("H43").Value = 21346710
$ws.Range("J43").Value = 49848484
$ws.Range("L43").Value = 49848484
$ws.Range("N43").Value = -49848870

$ws.Range("H68").Value = 2463
$ws.Range("I68").Value = 1262.25
$ws.Range("K68").Value = 1262.25
$ws.Range("M68").Value = -513.25

$ws.Range("H71").Value = 2463
$ws.Range("I71").Value = 1262.25
$ws.Range("K71").Value = 6311.25
$ws.Range("M71").Value = -2567.25

$ws.Range("H93").Value = 1374.6875
$ws.Range("I93").Value = 1276.174
$ws.Range("J93").Value = 1626.4445
$ws.Range("K93").Value = 1276.174
$ws.Range("L93").Value = 1626.4445
$ws.Range("M93").Value = -28.17399999999998
$ws.Range("N93").Value = -4122.4445

$ws.Range("H108").Value = 45842
$ws.Range("J108").Value = 45842
$ws.Range("L108").Value = 45842
$ws.Range("N108").Value = -53522

$ws.Range("H126").Value = 3666.3333
$ws.Range("I126").Value = 3249.5
$ws.Range("K126").Value = 9748.5
$ws.Range("M126").Value = -7278.5

$ws.Range("H132").Value = 8237.875
$ws.Range("I132").Value = 8627.066000000001
$ws.Range("K132").Value = 25881.198
$ws.Range("M132").Value = -23351.198

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 321032.25
$ws.Range("I14").Value = 1111.1111
$ws.Range("J14").Value = 2480500
$ws.Range("K14").Value = 1111.1111
$ws.Range("L14").Value = 2480500
$ws.Range("M14").Value = -943.1111000000001
$ws.Range("N14").Value = -2480836

$ws.Range("H113").Value = 2214.682
$ws.Range("I113").Value = 1067.2941
$ws.Range("J113").Value = 6115.8
$ws.Range("K113").Value = 3201.8823
$ws.Range("L113").Value = 18347.4
$ws.Range("M113").Value = -1031.8823
$ws.Range("N113").Value = -22687.4

$ws.Range("H122").Value = 3645.7
$ws.Range("I122").Value = 3614.4375
$ws.Range("K122").Value = 10843.3125
$ws.Range("M122").Value = -8393.3125

$ws.Range("H136").Value = 6430.9565
$ws.Range("I136").Value = 6971.1875
$ws.Range("J136").Value = 5196.143
$ws.Range("K136").Value = 20913.5625
$ws.Range("L136").Value = 15588.429
$ws.Range("M136").Value = -18363.5625
$ws.Range("N136").Value = -20688.429

$ws.Range("H140").Value = 89428
$ws.Range("J140").Value = 89428
$ws.Range("L140").Value = 89428
$ws.Range("N140").Value = -99788
